$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing 2017/2018/2019 header text from row 2 (C2:E2) down to row 3 (C3:E3).
$ws.Range("C2:E2").Copy()
$ws.Range("C3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Fill row 2 (B2:E2) with the new city name "Майкоп".
$ws.Range("B2:E2").Value = "Майкоп"

# Restore the view: no frozen top-left scroll offset, selection on H4.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H4").Select()
